$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet previously had two separate reward-coefficient parameter rows
# (wait_veh_reward_coef / wait_ped_reward_coef). They are being collapsed
# into a single "wait_reward_coef" parameter row, fixing a key error.
#
# Row 7 currently holds "wait_veh_reward_coef" and row 8 holds
# "wait_ped_reward_coef". We delete row 8 (shifting every row below it up
# by one) and rewrite row 7 so it becomes the new, unified parameter.

$ws.Rows("8").Delete()

$ws.Range("B7").Value = "wait_reward_coef"
$ws.Range("C7").Value = "Balancing coefficient for wait in reward calculation"
$ws.Range("D7").Value = 1

# The old C7 used a special rich-text/Consolas style (wrapped, taller row)
# to show the "alpha" symbol prefix. The merged description is now plain
# text, matching the style used by the rest of the parameter rows, so turn
# off wrapping and let the row height return to the sheet default.
$ws.Range("C7").WrapText = $false
$ws.Rows("7").EntireRow.AutoFit()

# Update the active selection left behind by the editor.
$ws.Range("E4").Select() | Out-Null
